$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the stray A1 label (the old "x" placeholder chord is being
#     promoted into a real row/column, so A1's corner cell becomes blank). ---
$ws.Range("A1").ClearContents()

# --- Rename the placeholder chord ("x" / "C/E") to the new "Cmaj7" chord,
#     used both as a column header (I1) and a row label (A9). ---
$ws.Range("I1").Value = "Cmaj7"
$ws.Range("A9").Value = "Cmaj7"

# --- Updated Markov transition counts (melody object data) ---
$ws.Range("B2").Value = 4
$ws.Range("E2").Value = 23
$ws.Range("G2").Value = 14
$ws.Range("I2").Value = 2

$ws.Range("C3").Value = 4
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 19
$ws.Range("I3").Value = 2

$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 4

$ws.Range("E5").Value = 4
$ws.Range("G5").Value = 12

$ws.Range("D6").Value = 5
$ws.Range("F6").Value = 4

$ws.Range("D7").Value = 7
$ws.Range("G7").Value = 4

$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 3

# Row 9 (Cmaj7) is fully rewritten with new data
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 11
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 12
$ws.Range("L9").Value = 13
$ws.Range("M9").Value = 22

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 15
$ws.Range("M10").Value = 16

$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 24
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 14
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = 14
$ws.Range("M11").Value = 8

$ws.Range("I12").Value = 2
$ws.Range("L12").Value = 4

$ws.Range("I13").Value = 3
$ws.Range("M13").Value = 4

# --- New "Sum" column (N) totalling each row's transition counts ---
$ws.Range("N1").Value = "Sum"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").Font.Underline = $true
$ws.Range("N2").Formula = "=SUM(B2:M2)"
$ws.Range("N3:N13").Formula = "=SUM(B3:M3)"

# --- Restore the active selection left behind by the author ---
$ws.Range("J17").Select()
